# "updated validations in add deals test cases"
#
# On the "deals" sheet (sheet3.xml):
#  - add two trailing columns: predictedCloseDate (O) / actualCloseDate (P),
#    filled with date values on rows 2-3, formatted as dates
#  - re-purpose the "probability" column (E) as free-text (values stored as
#    strings "80"/"60" instead of numbers), header + data cells get a
#    text number-format
#  - selection moves to column F

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new columns O (predictedCloseDate) / P (actualCloseDate) --------------

$ws.Range("O1").Value = "predictedCloseDate"
$ws.Range("O1").Interior.Color = 65535          # match the other header cells' yellow fill
$ws.Range("P1").Value = "actualCloseDate"
$ws.Range("P1").Interior.Color = 65535

$ws.Range("O2").Value = 43432                   # 2018-11-28
$ws.Range("O2").NumberFormat = "d-mmm-yy"
$ws.Range("P2").Value = 43434                   # 2018-11-30
$ws.Range("P2").NumberFormat = "d-mmm-yy"

$ws.Range("O3").Value = 43429                   # 2018-11-25
$ws.Range("O3").NumberFormat = "d-mmm-yy"
$ws.Range("P3").Value = 43430                   # 2018-11-26
$ws.Range("P3").NumberFormat = "d-mmm-yy"

$ws.Columns.Item(15).ColumnWidth = 18
$ws.Columns.Item(16).ColumnWidth = 14.6

# --- probability column (E) switches to free-text validation ---------------

$ws.Range("E1").NumberFormat = "@"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "80"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "60"

# --- selection moves onto column F ------------------------------------------

[void]$ws.Columns.Item(6).Select()
